$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 62529252
$ws.Range("I62").Value = 62529252
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 62529252
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -62528628
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 62529252
$ws.Range("I65").Value = 62529252
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 312646260
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -312643140
$ws.Range("N65").ClearContents()
$ws.Range("H100").Value = 12215.1
$ws.Range("I100").Value = 23689.555
$ws.Range("K100").Value = 23689.555
$ws.Range("M100").Value = -23148.555
$ws.Range("H135").Value = 1107.0526
$ws.Range("I135").Value = 1135.6111
$ws.Range("J135").Value = 593
$ws.Range("K135").Value = 10220.4999
$ws.Range("L135").Value = 5337
$ws.Range("M135").Value = -7685.499900000001
$ws.Range("N135").Value = -10407
$ws.Range("H137").Value = 1383.5122
$ws.Range("I137").Value = 1127.871
$ws.Range("J137").Value = 2176
$ws.Range("K137").Value = 3383.613
$ws.Range("L137").Value = 6528
$ws.Range("M137").Value = -833.6130000000003
$ws.Range("N137").Value = -11628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1584.3276
$ws.Range("I74").Value = 1701.878
$ws.Range("K74").Value = 1701.878
$ws.Range("M74").Value = -827.8779999999999
$ws.Range("H77").Value = 1584.3276
$ws.Range("I77").Value = 1701.878
$ws.Range("K77").Value = 8509.389999999999
$ws.Range("M77").Value = -4141.389999999999
$ws.Range("H86").Value = 514900
$ws.Range("J86").Value = 514900
$ws.Range("L86").Value = 514900
$ws.Range("N86").Value = -517272
$ws.Range("H89").Value = 514900
$ws.Range("J89").Value = 514900
$ws.Range("L89").Value = 1544700
$ws.Range("N89").Value = -1556556
$ws.Range("H97").Value = 1372.3158
$ws.Range("I97").Value = 1179.4
$ws.Range("J97").Value = 1586.6666
$ws.Range("K97").Value = 1179.4
$ws.Range("L97").Value = 1586.6666
$ws.Range("M97").Value = -683.4000000000001
$ws.Range("N97").Value = -2578.6666
$ws.Range("H132").Value = 17244012
$ws.Range("I132").Value = 20000682
$ws.Range("J132").Value = 14832.375
$ws.Range("K132").Value = 60002046
$ws.Range("L132").Value = 44497.125
$ws.Range("M132").Value = -59999516
$ws.Range("N132").Value = -49557.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 27029038
$ws.Range("I105").Value = 1981.6562
$ws.Range("K105").Value = 1981.6562
$ws.Range("M105").Value = -234.6561999999999
$ws.Range("H134").Value = 2370560.5
$ws.Range("I134").Value = 5716.185
$ws.Range("J134").Value = 5563100.5
$ws.Range("K134").Value = 17148.555
$ws.Range("L134").Value = 16689301.5
$ws.Range("M134").Value = -14613.555
$ws.Range("N134").Value = -16694371.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 41667800
$ws.Range("I58").Value = 62500868
$ws.Range("J58").Value = 1662.625
$ws.Range("K58").Value = 62500868
$ws.Range("L58").Value = 1662.625
$ws.Range("M58").Value = -62500665
$ws.Range("N58").Value = -2068.625
$ws.Range("H105").Value = 2481.6882
$ws.Range("I105").Value = 2395.4666
$ws.Range("J105").Value = 5715
$ws.Range("K105").Value = 2395.4666
$ws.Range("L105").Value = 5715
$ws.Range("M105").Value = -648.4666000000002
$ws.Range("N105").Value = -9209
$ws.Range("H134").Value = 1650.2646
$ws.Range("I134").Value = 1675.2812
$ws.Range("K134").Value = 5025.8436
$ws.Range("M134").Value = -2490.8436
$ws.Range("H136").Value = 41667800
$ws.Range("I136").Value = 62500868
$ws.Range("J136").Value = 1662.625
$ws.Range("K136").Value = 187502604
$ws.Range("L136").Value = 4987.875
$ws.Range("M136").Value = -187500054
$ws.Range("N136").Value = -10087.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 93
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 119
$ws.Range("K23").Value = 45
$ws.Range("L23").Value = 357
$ws.Range("M23").Value = 190
$ws.Range("N23").Value = -827
$ws.Range("H34").Value = 1125
$ws.Range("I34").Value = 220
$ws.Range("J34").Value = 1527.2222
$ws.Range("K34").Value = 660
$ws.Range("L34").Value = 4581.6666
$ws.Range("M34").Value = -576
$ws.Range("N34").Value = -4749.6666
$ws.Range("H39").Value = 1900
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1900
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5700
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6288
$ws.Range("H55").Value = 1786.4286
$ws.Range("J55").Value = 1786.4286
$ws.Range("L55").Value = 5359.2858
$ws.Range("N55").Value = -5713.2858
$ws.Range("H64").Value = 1388.4
$ws.Range("J64").Value = 1497.25
$ws.Range("L64").Value = 4491.75
$ws.Range("N64").Value = -5031.75
$ws.Range("H67").Value = 1388.4
$ws.Range("J67").Value = 1497.25
$ws.Range("L67").Value = 4491.75
$ws.Range("N67").Value = -6363.75
$ws.Range("H130").Value = 83334970
$ws.Range("I130").Value = 500000000
$ws.Range("K130").Value = 1500000000
$ws.Range("M130").Value = -1499994980
$ws.Range("H131").Value = 752.4545000000001
$ws.Range("I131").Value = 464.44446
$ws.Range("J131").Value = 781.25555
$ws.Range("K131").Value = 1393.33338
$ws.Range("L131").Value = 2343.76665
$ws.Range("M131").Value = 3646.66662
$ws.Range("N131").Value = -12423.76665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3907.9756
$ws.Range("I70").Value = 3867.1072
$ws.Range("J70").Value = 3996
$ws.Range("K70").Value = 3867.1072
$ws.Range("L70").Value = 3996
$ws.Range("M70").Value = -3597.1072
$ws.Range("N70").Value = -4536
$ws.Range("H73").Value = 3907.9756
$ws.Range("I73").Value = 3867.1072
$ws.Range("J73").Value = 3996
$ws.Range("K73").Value = 3867.1072
$ws.Range("L73").Value = 3996
$ws.Range("M73").Value = -2931.1072
$ws.Range("N73").Value = -5868
$ws.Range("H108").Value = 35110.668
$ws.Range("J108").Value = 35110.668
$ws.Range("L108").Value = 35110.668
$ws.Range("N108").Value = -42790.668
$ws.Range("H126").Value = 4412.933
$ws.Range("I126").Value = 3639.8
$ws.Range("K126").Value = 10919.4
$ws.Range("M126").Value = -8449.400000000001
$ws.Range("H132").Value = 5618.394
$ws.Range("I132").Value = 1228.7916
$ws.Range("J132").Value = 17324
$ws.Range("K132").Value = 3686.3748
$ws.Range("L132").Value = 51972
$ws.Range("M132").Value = -1156.3748
$ws.Range("N132").Value = -57032

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4223009
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590
$ws.Range("H27").Value = 4223009
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214
$ws.Range("H40").Value = 125001000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H93").Value = 1343.7646
$ws.Range("I93").Value = 1341.9231
$ws.Range("J93").Value = 1349.75
$ws.Range("K93").Value = 1341.9231
$ws.Range("L93").Value = 1349.75
$ws.Range("M93").Value = -93.92309999999998
$ws.Range("N93").Value = -3845.75
$ws.Range("H100").Value = 3936.0588
$ws.Range("I100").Value = 4000.7144
$ws.Range("K100").Value = 4000.7144
$ws.Range("M100").Value = -3459.7144
$ws.Range("H108").Value = 33296
$ws.Range("J108").Value = 33296
$ws.Range("L108").Value = 33296
$ws.Range("N108").Value = -40976
$ws.Range("H122").Value = 44999.832
$ws.Range("I122").Value = 51999.8
$ws.Range("K122").Value = 155999.4
$ws.Range("M122").Value = -153549.4
$ws.Range("H132").Value = 14710861
$ws.Range("I132").Value = 30305162
$ws.Range("J132").Value = 7663.086
$ws.Range("K132").Value = 90915486
$ws.Range("L132").Value = 22989.258
$ws.Range("M132").Value = -90912956
$ws.Range("N132").Value = -28049.258
$ws.Range("H136").Value = 5867.7144
$ws.Range("I136").Value = 4289.8335
$ws.Range("K136").Value = 12869.5005
$ws.Range("M136").Value = -10319.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 62500856
$ws.Range("I81").Value = 62500856
$ws.Range("K81").Value = 125001712
$ws.Range("M81").Value = -125000651
$ws.Range("H84").Value = 62500856
$ws.Range("I84").Value = 62500856
$ws.Range("K84").Value = 625008560
$ws.Range("M84").Value = -625003256
$ws.Range("H132").Value = 21765480
$ws.Range("I132").Value = 38502600
$ws.Range("K132").Value = 115507800
$ws.Range("M132").Value = -115505270
$ws.Range("H136").Value = 5199.604
$ws.Range("I136").Value = 11279
$ws.Range("J136").Value = 1551.9667
$ws.Range("K136").Value = 33837
$ws.Range("L136").Value = 4655.9001
$ws.Range("M136").Value = -31287
$ws.Range("N136").Value = -9755.900099999999
